$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-06 Monday" "2025-01-07 Tuesday"

Replace-Text "605÷4=" "293÷8="
Replace-Text "178÷3=" "133÷2="
Replace-Text "759÷7=" "566÷9="
Replace-Text "751÷9=" "811÷3="
Replace-Text "198÷2=" "295÷7="
Replace-Text "234÷9=" "361÷6="
Replace-Text "220÷6=" "102÷9="
Replace-Text "355÷4=" "269÷4="
Replace-Text "842÷2=" "266÷6="
Replace-Text "514÷2=" "773÷9="
Replace-Text "728÷9=" "626÷2="
Replace-Text "248÷8=" "355÷3="
Replace-Text "986÷2=" "138÷9="
Replace-Text "250÷2=" "112÷5="
Replace-Text "768÷8=" "985÷3="
Replace-Text "115÷7=" "802÷7="
Replace-Text "129÷2=" "661÷2="
Replace-Text "321÷2=" "548÷3="
Replace-Text "398÷2=" "446÷8="
Replace-Text "569÷4=" "659÷6="
Replace-Text "588÷2=" "369÷6="
Replace-Text "758÷2=" "849÷4="
Replace-Text "772÷6=" "108÷3="
Replace-Text "969÷3=" "687÷8="
Replace-Text "174÷7=" "418÷8="
